# Updates the "cryptos" price list with refreshed prices / 1h volume
# percentages (and, for rows 44-45, the two coins that swapped rank).
#
# All of the target cells are stored as plain text in the workbook (coin
# prices like "59.469.70" or "0.999" are not valid Excel numbers/are meant
# to stay text), so each write briefly forces the cell to Text format,
# assigns the literal string, then restores "Normal" style so no stray
# per-cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "59.469.70"
Set-TextValue "E2" "  +0.12%  "

Set-TextValue "D3" "2.636.99"
Set-TextValue "E3" "  +1.15%  "

Set-TextValue "E4" "  +0.11%  "

Set-TextValue "D5" "535.99"
Set-TextValue "E5" "  -0.19%  "

Set-TextValue "D6" "144.81"
Set-TextValue "E6" "  +2.71%  "

Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  +0.08%  "

Set-TextValue "D8" "0.572"
Set-TextValue "E8" "  +0.53%  "

Set-TextValue "D9" "2.646.24"
Set-TextValue "E9" "  +1.07%  "

Set-TextValue "D10" "6.62"
Set-TextValue "E10" "  +2.32%  "

Set-TextValue "D11" "0.102"
Set-TextValue "E11" "  -0.87%  "

Set-TextValue "D12" "0.337"
Set-TextValue "E12" "  -0.17%  "

Set-TextValue "E13" "  -0.69%  "

Set-TextValue "D14" "3.112.24"
Set-TextValue "E14" "  +1.45%  "

Set-TextValue "D15" "59.396.86"
Set-TextValue "E15" "  +0.11%  "

Set-TextValue "D16" "21.07"
Set-TextValue "E16" "  +2.46%  "

Set-TextValue "D17" "2.630.08"
Set-TextValue "E17" "  +0.63%  "

Set-TextValue "D18" "0.0000134"
Set-TextValue "E18" "  +0.29%  "

Set-TextValue "D19" "339.95"
Set-TextValue "E19" "  -1.72%  "

Set-TextValue "D20" "4.39"
Set-TextValue "E20" "  +0.75%  "

Set-TextValue "D21" "10.37"
Set-TextValue "E21" "  +2.17%  "

Set-TextValue "D22" "6.28"
Set-TextValue "E22" "  -1.87%  "

Set-TextValue "D24" "66.94"
Set-TextValue "E24" "  -0.22%  "

Set-TextValue "D25" "0.414"
Set-TextValue "E25" "  +1.29%  "

Set-TextValue "D26" "0.165"
Set-TextValue "E26" "  -1.42%  "

Set-TextValue "D28" "7.28"
Set-TextValue "E28" "  +1.01%  "

Set-TextValue "D29" "0.0₃0745"
Set-TextValue "E29" "  -0.43%  "

Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  +0.00%  "

Set-TextValue "E31" "  +0.62%  "

Set-TextValue "D32" "5.83"
Set-TextValue "E32" "  -0.69%  "

Set-TextValue "D33" "18.87"
Set-TextValue "E33" "  -0.18%  "

Set-TextValue "D34" "150.76"
Set-TextValue "E34" "  +1.16%  "

Set-TextValue "D35" "3.99"
Set-TextValue "E35" "  -0.25%  "

Set-TextValue "D36" "1.13"
Set-TextValue "E36" "  +0.91%  "

Set-TextValue "D37" "0.834"
Set-TextValue "E37" "  -0.92%  "

Set-TextValue "D38" "0.835"
Set-TextValue "E38" "  -0.59%  "

Set-TextValue "D39" "1.45"
Set-TextValue "E39" "  -0.98%  "

Set-TextValue "D40" "289.57"
Set-TextValue "E40" "  +4.65%  "

Set-TextValue "D41" "3.59"
Set-TextValue "E41" "  +0.88%  "

Set-TextValue "E42" "  +0.14%  "

Set-TextValue "D43" "0.603"
Set-TextValue "E43" "  +0.51%  "

Set-TextValue "B44" "WhiteBITCoin"
Set-TextValue "C44" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D44" "10.73"
Set-TextValue "E44" "  -0.26%  "

Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "19.32"
Set-TextValue "E45" "  +3.45%  "

Set-TextValue "D46" "0.0535"
Set-TextValue "E46" "  +1.96%  "

Set-TextValue "D47" "0.0946"
Set-TextValue "E47" "  -1.73%  "

Set-TextValue "D48" "1.971.97"
Set-TextValue "E48" "  +1.20%  "

Set-TextValue "D49" "0.0226"
Set-TextValue "E49" "  +1.03%  "

Set-TextValue "D50" "4.54"
Set-TextValue "E50" "  +0.21%  "

Set-TextValue "D51" "18.30"
Set-TextValue "E51" "  -0.34%  "
